$wb = $excel.ActiveWorkbook

# --- Sheet "Home win" ---
$ws = $wb.Worksheets.Item("Home win")
$ws.Range("A2").Value = "06-01-2025 20:00"
$ws.Range("B2").Value = "ENGLAND"
$ws.Range("C2").Value = "CHAMPIONSHIP"
$ws.Range("D2").Value = "QPR - Luton"
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 2.4

$ws.Range("A3").Value = "06-01-2025 22:30"
$ws.Range("B3").Value = "BRAZIL"
$ws.Range("C3").Value = "SÃO PAULO YOUTH CUP"
$ws.Range("D3").Value = "Jaciobá U20 - Santos U20"
$ws.Range("E3").Value = 70
$ws.Range("F3").Value = 46

$ws.Range("A4").Value = "06-01-2025 12:30"
$ws.Range("B4").Value = "SAUDI-ARABIA"
$ws.Range("C4").Value = "DIVISION 1"
$ws.Range("D4").Value = "Al Najma - Abha"
$ws.Range("E4").Value = 80
$ws.Range("F4").Value = 2.05

$ws.Range("A5").Value = "06-01-2025 14:40"
$ws.Range("B5").Value = "SAUDI-ARABIA"
$ws.Range("C5").Value = "DIVISION 1"
$ws.Range("D5").Value = "Al Jubail - Ohod"
$ws.Range("E5").Value = 70
$ws.Range("F5").Value = 1.91

# --- Sheet "Draw" ---
$ws = $wb.Worksheets.Item("Draw")
$ws.Range("A2").Value = "07-01-2025 19:45"
$ws.Range("B2").Value = "ENGLAND"
$ws.Range("C2").Value = "LEAGUE ONE"
$ws.Range("D2").Value = "Wycombe - Huddersfield"
$ws.Range("E2").Value = 60
$ws.Range("F2").Value = 3.5

$ws.Range("A3:F6").Clear()

# --- Sheet "Btts" ---
$ws = $wb.Worksheets.Item("Btts")
$ws.Range("A2").Value = "06-01-2025 21:15"
$ws.Range("B2").Value = "BRAZIL"
$ws.Range("C2").Value = "SÃO PAULO YOUTH CUP"
$ws.Range("D2").Value = "Juventude U20 - América RN U20"
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 1.85

$ws.Range("A3").Value = "06-01-2025 14:40"
$ws.Range("B3").Value = "SAUDI-ARABIA"
$ws.Range("C3").Value = "DIVISION 1"
$ws.Range("D3").Value = "Al Jubail - Ohod"
$ws.Range("E3").Value = 80
$ws.Range("F3").Value = 1.91

$ws.Range("A4:F11").Clear()

# --- Sheet "Away Win" ---
$ws = $wb.Worksheets.Item("Away Win")
$ws.Range("A2").Value = "07-01-2025 19:45"
$ws.Range("B2").Value = "NORTHERN-IRELAND"
$ws.Range("C2").Value = "LEAGUE CUP"
$ws.Range("D2").Value = "Ballymena United - Glentoran"
$ws.Range("E2").Value = 73.3
$ws.Range("F2").Value = 1.95

$ws.Range("A3:F4").Clear()
